$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Estadisticos 2P" (segundo parcial results)
# ---------------------------------------------------------------------------
$ws2p = $wb.Worksheets.Item("Estadisticos 2P")
$ws2p.Cells.Item(2, 4).Value = 0
$ws2p.Cells.Item(2, 5).Value = 10
$ws2p.Cells.Item(2, 6).Value = 14
$ws2p.Cells.Item(2, 7).Value = 58.33
$ws2p.Cells.Item(2, 8).Value = 6.7

$ws2p.Cells.Item(3, 4).Value = 0
$ws2p.Cells.Item(3, 5).Value = 9
$ws2p.Cells.Item(3, 6).Value = 22
$ws2p.Cells.Item(3, 7).Value = 70.97
$ws2p.Cells.Item(3, 8).Value = 7.4

$ws2p.Cells.Item(4, 4).Value = 0
$ws2p.Cells.Item(4, 5).Value = 12
$ws2p.Cells.Item(4, 6).Value = 8
$ws2p.Cells.Item(4, 7).Value = 40
$ws2p.Cells.Item(4, 8).Value = 5.9

# ---------------------------------------------------------------------------
# Sheet "Estadisticos Final" (combined / final results)
# ---------------------------------------------------------------------------
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")
$wsFinal.Cells.Item(2, 5).Value = 10
$wsFinal.Cells.Item(2, 6).Value = 14
$wsFinal.Cells.Item(2, 7).Value = 58.33
$wsFinal.Cells.Item(2, 8).Value = 6.5

$wsFinal.Cells.Item(3, 5).Value = 9
$wsFinal.Cells.Item(3, 6).Value = 22
$wsFinal.Cells.Item(3, 7).Value = 70.97
$wsFinal.Cells.Item(3, 8).Value = 7.3

$wsFinal.Cells.Item(4, 5).Value = 12
$wsFinal.Cells.Item(4, 6).Value = 8
$wsFinal.Cells.Item(4, 7).Value = 40
$wsFinal.Cells.Item(4, 8).Value = 6

# ---------------------------------------------------------------------------
# Sheet "Rescatables" (list of students eligible for remedial exam, refreshed
# after the second partial exam)
# ---------------------------------------------------------------------------
$wsResc = $wb.Worksheets.Item("Rescatables")

$data = @(
  @(21330051920007, 'COBOS', 'NOLASCO', 'YOLET', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6AEV', 4),
  @(22330051920413, 'LOBATO', 'ANTONIO', 'FABIAN ALEJANDRO', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6AEV', 4),
  @(23330051920045, 'SANTIAGO', 'GARCIA', 'URIEL', 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO', '4BEM', 3),
  @(22330051920177, 'CAMPOS', 'CABRERA', 'MARCO', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6AEV', 3),
  @(22330051920188, 'HERNANDEZ', 'ROJAS', 'DAVID', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6AEV', 3),
  @(22330051920190, 'JIMENEZ', 'CIRUELO', 'ARACELY', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6AEV', 3),
  @(22330051920371, 'VERA', 'GONZALEZ', 'ISRAEL', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6AEV', 3),
  @(22330051920031, 'CASTILLO', 'GONZALEZ', 'RICARDO', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6BEM', 3),
  @(22330051920038, 'HERRERA', 'ACOSTA', 'MIGUEL ANTONIO', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6BEM', 3),
  @(22330051920189, 'JENKINS', 'GARCIA', 'ARTHUR RICHARD', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6BEM', 3),
  @(22330051920043, 'PALOMINO', 'HERNANDEZ', 'AARON MIGUEL', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6BEM', 3),
  @(23330051920028, 'CAMPOS', 'RIVERA', 'IRVING', 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO', '4BEM', 2),
  @(22330051920007, 'CARRERA', 'GARCIA', 'ANA KAREN', 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO', '4BEM', 2),
  @(23330051920037, 'HERNANDEZ', 'MARCELINO', 'LEONEL', 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO', '4BEM', 2),
  @(22330051920021, 'MEJIA', 'CRUZ', 'JOSE FRANCISCO', 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO', '4BEM', 2),
  @(23330051920212, 'VERA', 'VILLA', 'ALEX URIEL', 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO', '4BEM', 2),
  @(22330051920359, 'MARQUEZ', 'TIZA', 'CRISTIAN OSMAR', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6AEV', 2),
  @(22330051920192, 'MENDOZA', 'HERNANDEZ', 'ERIK OMAR', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6AEV', 2),
  @(22330051920356, 'RUIZ', 'ALFONSO', 'JOSUE GUSTAVO', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6AEV', 2),
  @(22330051920033, 'CRESCENCIO', 'DIAZ', 'DIEGO ARMANDO', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6BEM', 2),
  @(22330051920191, 'MELCHOR', 'DE LA CRUZ', 'ALDAHIR', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6BEM', 2),
  @(22330051920327, 'PEREZ', 'CONTRERAS', 'JORGE IVAN', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6BEM', 2),
  @(23330051920032, 'DE JESUS', 'VERA', 'EDUARDO', 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO', '4BEM', 1),
  @(23330051920040, 'MENDEZ', 'SARMIENTO', 'ALAN URIEL', 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO', '4BEM', 1),
  @(23330051920043, 'RODRIGUEZ', 'SOLANO', 'IAN', 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO', '4BEM', 1),
  @(23330051920046, 'TRUJILLO', 'CALIHUA', 'YAEL ISSAI', 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO', '4BEM', 1),
  @(22330051920180, 'DE LOS SANTOS', 'HERNANDEZ', 'ABDIEL NOE', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6AEV', 1),
  @(22330051920201, 'XOTLANIHUA', 'COLOHUA', 'ALEXANDER', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6AEV', 1),
  @(22330051920202, 'XOTLANIHUA', 'COLOHUA', 'ERIK', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6AEV', 1),
  @(22330051920034, 'CRUZ', 'REYES', 'CARLOS YAEL', 'REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS', '6BEM', 1)

)

$r = 2
foreach ($row in $data) {
    $wsResc.Cells.Item($r, 1).Value = $row[0]
    $wsResc.Cells.Item($r, 2).Value = $row[1]
    $wsResc.Cells.Item($r, 3).Value = $row[2]
    $wsResc.Cells.Item($r, 4).Value = $row[3]
    $wsResc.Cells.Item($r, 5).Value = $row[4]
    $wsResc.Cells.Item($r, 6).Value = $row[5]
    $wsResc.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}
